$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.72   # Current Capital
$summary.Range("B4").Value = -0.28     # Total P&L $
$summary.Range("B5").Value = -0.2      # Total P&L %
$summary.Range("B6").Value = 28        # Total Trades
$summary.Range("B8").Value = 14        # Losing Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.72      # Capital
$status.Range("D4").Value = 28         # Trades
$status.Range("E4").Value = -0.28      # P&L $
$status.Range("F4").Value = -0.28      # P&L %
$status.Range("G4").Value = 28.57      # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade row (#28) to a trades sheet
# ---------------------------------------------------------------------------
function Add-TradeRow($sheet) {
    $sheet.Range("A29").Value = 28
    $sheet.Range("B29").Value = "'2026-02-17"
    $sheet.Range("C29").Value = "15:22:07"
    $sheet.Range("D29").Value = "MarketMaking"
    $sheet.Range("E29").Value = "UP"
    $sheet.Range("F29").Value = 0.79
    $sheet.Range("G29").Value = 0.65
    $sheet.Range("H29").Value = "CLOSED"
    $sheet.Range("I29").Value = -17.7215
    $sheet.Range("J29").Value = -0.14
    $sheet.Range("K29").Value = 99.72
    $sheet.Range("L29").Value = 0
    $sheet.Range("M29").Value = 0
    $sheet.Range("N29").Value = 0.6
    $sheet.Range("O29").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P29").Value = "early_exit"
    $sheet.Range("Q29").Value = 0.12
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
